$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$lo = $ws.ListObjects.Item(1)

# Rows that were previously hidden and now become visible, gaining a
# tax_code value of "LR022" in column E.
$rowsToReveal = @(4, 5, 7, 8, 9, 14, 17, 18, 20, 21, 22, 27, 30)

foreach ($r in $rowsToReveal) {
    $ws.Rows.Item($r).Hidden = $false
    $ws.Cells.Item($r, 5).Value = "LR022"
}

# Rows 19 and 23 already had a tax_code value (LR003); update it to LR022.
$ws.Cells.Item(19, 5).Value = "LR022"
$ws.Cells.Item(23, 5).Value = "LR022"

# Clear the autofilter criteria on the "property_class" column (field 2)
# while keeping the table's autofilter dropdowns in place.
$lo.Range.AutoFilter(2, @())

# Update the selected cell in the sheet view from F1 to E30.
[void]$ws.Range("E30").Select()
